# Apply revised implementation of link process partitioning values
# (eGRID emission factor data for Cl) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 (AZNM)
$ws.Range("B7").Value = 2060.9468865944887
$ws.Range("C7").Value = 73592.327179765198
$ws.Range("D7").Value = 1840.4857955840102

# Row 11 (MROE)
$ws.Range("B11").Value = 1029.8374098269771
$ws.Range("D11").Value = 59490.093710174246

# Row 12 (MROW)
$ws.Range("B12").Value = 1118.6298522398338
$ws.Range("D12").Value = 52675.340891975153

# Row 14 (NWPP)
$ws.Range("B14").Value = 1030.1160733044671
$ws.Range("D14").Value = 32270.382857352699

# Row 17 (RFCE)
$ws.Range("B17").Value = 2969.3100070159107
$ws.Range("C17").Value = 108858.86699265332
$ws.Range("D17").Value = 2534.6076892404935

# Row 18 (RFCM)
$ws.Range("B18").Value = 2930.6553436945082
$ws.Range("C18").Value = 107623.98025053306
$ws.Range("D18").Value = 5595.1542128661213

# Row 19 (RFCW)
$ws.Range("B19").Value = 1024.9619652441672
$ws.Range("D19").Value = 62332.306466173672

# Row 20 (RMPA)
$ws.Range("B20").Value = 6055.1228261945453
$ws.Range("C20").Value = 212559.3441949547
$ws.Range("D20").Value = 6937.5134681139371

# Row 21 (SPNO)
$ws.Range("B21").Value = 124.52712194958033

# Row 22 (SPSO)
$ws.Range("B22").Value = 1082.6159223638724
$ws.Range("D22").Value = 66782.766130014672

# Row 23 (SRMV)
$ws.Range("B23").Value = 1008.6125205841707
$ws.Range("D23").Value = 61171.743188536529

# Row 24 (SRMW)
$ws.Range("B24").Value = 1077.1660490929664
$ws.Range("D24").Value = 60068.167449748937

# Row 25 (SRSO)
$ws.Range("B25").Value = 638.14336637155009
$ws.Range("D25").Value = 59882.547272092561

# Row 26 (SRTV)
$ws.Range("B26").Value = 1717.3917210581315
$ws.Range("C26").Value = 55135.954922450692
$ws.Range("D26").Value = 2089.1864249921691

# Row 27 (SRVC)
$ws.Range("B27").Value = 1274.4336159711088
$ws.Range("D27").Value = 60418.339181479845

# Row 28 (next region)
$ws.Range("B28").Value = 6712.6367083421674
$ws.Range("C28").Value = 257725.97748798455
$ws.Range("D28").Value = 6404.8853861257539
